# Scheduled-runner update: refresh Universalis market-price snapshots and
# recompute Leve crafting-profit columns (H:N) for the affected job sheets.
# Generated from the authoritative cell-level diff; LTW has no changes this run.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 101.25
$ws.Range("I5").Value = 101.25
$ws.Range("K5").Value = 101.25
$ws.Range("M5").Value = 13.75
$ws.Range("H18").Value = 233.33333
$ws.Range("I18").Value = 233.33333
$ws.Range("K18").Value = 233.33333
$ws.Range("M18").Value = 50.66667000000001
$ws.Range("H64").Value = 6832.3335
$ws.Range("J64").Value = 6832.3335
$ws.Range("L64").Value = 6832.3335
$ws.Range("N64").Value = -7328.3335
$ws.Range("H67").Value = 6832.3335
$ws.Range("J67").Value = 6832.3335
$ws.Range("L67").Value = 6832.3335
$ws.Range("N67").Value = -8548.333500000001
$ws.Range("H82").Value = 433
$ws.Range("I82").Value = 433
$ws.Range("K82").Value = 1299
$ws.Range("M82").Value = -893
$ws.Range("H85").Value = 433
$ws.Range("I85").Value = 433
$ws.Range("K85").Value = 1299
$ws.Range("M85").Value = 105
$ws.Range("H115").Value = 572.25
$ws.Range("I115").Value = 572.25
$ws.Range("K115").Value = 1716.75
$ws.Range("M115").Value = -149.75
$ws.Range("H137").Value = 6113.321
$ws.Range("I137").Value = 2258.8809
$ws.Range("J137").Value = 20830.273
$ws.Range("K137").Value = 6776.6427
$ws.Range("L137").Value = 62490.819
$ws.Range("M137").Value = -4226.6427
$ws.Range("N137").Value = -67590.819

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 354.375
$ws.Range("I4").Value = 305
$ws.Range("J4").Value = 700
$ws.Range("K4").Value = 305
$ws.Range("L4").Value = 700
$ws.Range("M4").Value = -189
$ws.Range("N4").Value = -932
$ws.Range("H45").Value = 4460.25
$ws.Range("I45").Value = 4447
$ws.Range("K45").Value = 4447
$ws.Range("M45").Value = -4070
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H61").Value = 707860.1
$ws.Range("I61").Value = 2887.4075
$ws.Range("J61").Value = 1765319.1
$ws.Range("K61").Value = 2887.4075
$ws.Range("L61").Value = 1765319.1
$ws.Range("M61").Value = -2675.4075
$ws.Range("N61").Value = -1765743.1
$ws.Range("H136").Value = 707860.1
$ws.Range("I136").Value = 2887.4075
$ws.Range("J136").Value = 1765319.1
$ws.Range("K136").Value = 8662.2225
$ws.Range("L136").Value = 5295957.300000001
$ws.Range("M136").Value = -6112.2225
$ws.Range("N136").Value = -5301057.300000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 181.35715
$ws.Range("I7").Value = 263.1111
$ws.Range("J7").Value = 34.2
$ws.Range("K7").Value = 263.1111
$ws.Range("L7").Value = 34.2
$ws.Range("M7").Value = -150.1111
$ws.Range("N7").Value = -260.2
$ws.Range("H11").Value = 1545.3846
$ws.Range("J11").Value = 2607.8572
$ws.Range("L11").Value = 2607.8572
$ws.Range("N11").Value = -2887.8572
$ws.Range("H12").Value = 1312.5
$ws.Range("J12").Value = 1733.3334
$ws.Range("L12").Value = 1733.3334
$ws.Range("N12").Value = -2069.3334
$ws.Range("H20").Value = 20272.143
$ws.Range("I20").Value = 6203.0557
$ws.Range("J20").Value = 30823.959
$ws.Range("K20").Value = 6203.0557
$ws.Range("L20").Value = 30823.959
$ws.Range("M20").Value = -5956.0557
$ws.Range("N20").Value = -31317.959
$ws.Range("H134").Value = 46959.223
$ws.Range("I134").Value = 56436.5
$ws.Range("J134").Value = 28004.666
$ws.Range("K134").Value = 169309.5
$ws.Range("L134").Value = 84013.99800000001
$ws.Range("M134").Value = -166774.5
$ws.Range("N134").Value = -89083.99800000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 839.5625
$ws.Range("I22").Value = 370.9091
$ws.Range("K22").Value = 370.9091
$ws.Range("M22").Value = -20.90910000000002
$ws.Range("H31").Value = 11398.333
$ws.Range("I31").Value = 826.0476
$ws.Range("K31").Value = 826.0476
$ws.Range("M31").Value = -531.0476
$ws.Range("H34").Value = 11398.333
$ws.Range("I34").Value = 826.0476
$ws.Range("K34").Value = 826.0476
$ws.Range("M34").Value = -624.0476
$ws.Range("H86").Value = 24798.4
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 24798.4
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H99").Value = 3109.6155
$ws.Range("I99").Value = 1848.2858
$ws.Range("K99").Value = 1848.2858
$ws.Range("M99").Value = -350.2858000000001
$ws.Range("H126").Value = 3109.6155
$ws.Range("I126").Value = 1848.2858
$ws.Range("K126").Value = 5544.857400000001
$ws.Range("M126").Value = -3074.857400000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 35000
$ws.Range("J74").Value = 35000
$ws.Range("L74").Value = 105000
$ws.Range("N74").Value = -107122
$ws.Range("H77").Value = 35000
$ws.Range("J77").Value = 35000
$ws.Range("L77").Value = 315000
$ws.Range("N77").Value = -325608
$ws.Range("H107").Value = 1149.6522
$ws.Range("I107").Value = 454.23077
$ws.Range("J107").Value = 2053.7
$ws.Range("K107").Value = 1362.69231
$ws.Range("L107").Value = 6161.099999999999
$ws.Range("M107").Value = 557.3076900000001
$ws.Range("N107").Value = -10001.1
$ws.Range("H113").Value = 1123.5294
$ws.Range("I113").Value = 1062.375
$ws.Range("K113").Value = 3187.125
$ws.Range("M113").Value = -1017.125
$ws.Range("H122").Value = 11958984
$ws.Range("I122").Value = 37373948
$ws.Range("J122").Value = 2183997.2
$ws.Range("K122").Value = 336365532
$ws.Range("L122").Value = 19655974.8
$ws.Range("M122").Value = -336363082
$ws.Range("N122").Value = -19660874.8
$ws.Range("H129").Value = 16668615
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 1476.44
$ws.Range("I131").Value = 1418
$ws.Range("J131").Value = 1478.875
$ws.Range("K131").Value = 4254
$ws.Range("L131").Value = 4436.625
$ws.Range("M131").Value = 786
$ws.Range("N131").Value = -14516.625

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1055.5834
$ws.Range("I3").Value = 358
$ws.Range("J3").Value = 1753.1666
$ws.Range("K3").Value = 358
$ws.Range("L3").Value = 1753.1666
$ws.Range("M3").Value = -242
$ws.Range("N3").Value = -1985.1666
$ws.Range("H13").Value = 950.9
$ws.Range("I13").Value = 118.333336
$ws.Range("J13").Value = 1307.7142
$ws.Range("K13").Value = 118.333336
$ws.Range("L13").Value = 1307.7142
$ws.Range("M13").Value = 20.666664
$ws.Range("N13").Value = -1585.7142
$ws.Range("H14").Value = 2655.0557
$ws.Range("I14").Value = 2524.4
$ws.Range("J14").Value = 3308.3333
$ws.Range("K14").Value = 2524.4
$ws.Range("L14").Value = 3308.3333
$ws.Range("M14").Value = -2356.4
$ws.Range("N14").Value = -3644.3333
$ws.Range("H22").Value = 2000
$ws.Range("J22").Value = 4250
$ws.Range("L22").Value = 4250
$ws.Range("N22").Value = -5308
$ws.Range("H70").Value = 4921.8335
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4921.8335
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 4921.8335
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -5461.8335
$ws.Range("H73").Value = 4921.8335
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4921.8335
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 4921.8335
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -6793.8335
$ws.Range("H80").Value = 10896.385
$ws.Range("J80").Value = 13379.6
$ws.Range("L80").Value = 13379.6
$ws.Range("N80").Value = -15375.6
$ws.Range("H83").Value = 10896.385
$ws.Range("J83").Value = 13379.6
$ws.Range("L83").Value = 66898
$ws.Range("N83").Value = -76882
$ws.Range("H102").Value = 6489.769
$ws.Range("I102").Value = 6760.727
$ws.Range("J102").Value = 4999.5
$ws.Range("K102").Value = 6760.727
$ws.Range("L102").Value = 4999.5
$ws.Range("M102").Value = -5138.727
$ws.Range("N102").Value = -8243.5
$ws.Range("H123").Value = 55210.445
$ws.Range("J123").Value = 55210.445
$ws.Range("L123").Value = 55210.445
$ws.Range("N123").Value = -60110.445
$ws.Range("H134").Value = 55428.57
$ws.Range("J134").Value = 55428.57
$ws.Range("L134").Value = 166285.71
$ws.Range("N134").Value = -171355.71

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1652.125
$ws.Range("I23").Value = 411.16666
$ws.Range("J23").Value = 5375
$ws.Range("K23").Value = 411.16666
$ws.Range("L23").Value = 5375
$ws.Range("M23").Value = -182.16666
$ws.Range("N23").Value = -5833
$ws.Range("H81").Value = 983.8333
$ws.Range("I81").Value = 967
$ws.Range("K81").Value = 1934
$ws.Range("M81").Value = -873
$ws.Range("H84").Value = 983.8333
$ws.Range("I84").Value = 967
$ws.Range("K84").Value = 9670
$ws.Range("M84").Value = -4366
$ws.Range("H96").Value = 1846.7142
$ws.Range("J96").Value = 1862.8334
$ws.Range("L96").Value = 1862.8334
$ws.Range("N96").Value = -4608.8334
$ws.Range("H122").Value = 3936.7
$ws.Range("I122").Value = 2242.6843
$ws.Range("K122").Value = 6728.0529
$ws.Range("M122").Value = -4278.0529
